$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 61
$ws.Range("H61").Value = 140
$ws.Range("I61").Value = 140
$ws.Range("K61").Value = 420
$ws.Range("M61").Value = -248
# Row 136
$ws.Range("H136").Value = 46413.715
$ws.Range("J136").Value = 46413.715
$ws.Range("L136").Value = 46413.715
$ws.Range("N136").Value = -56613.715
# Row 137
$ws.Range("H137").Value = 1671051.1
$ws.Range("I137").Value = 4350379
$ws.Range("J137").Value = 5522.973
$ws.Range("K137").Value = 13051137
$ws.Range("L137").Value = 16568.919
$ws.Range("M137").Value = -13048587
$ws.Range("N137").Value = -21668.919

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1312.6818
$ws.Range("I2").Value = 938.8163500000001
$ws.Range("J2").Value = 2390.2942
$ws.Range("K2").Value = 938.8163500000001
$ws.Range("L2").Value = 2390.2942
$ws.Range("M2").Value = -825.8163500000001
$ws.Range("N2").Value = -2616.2942
# Row 63
$ws.Range("H63").Value = 2785
$ws.Range("I63").Value = 2830.5
$ws.Range("J63").Value = 2633.3333
$ws.Range("K63").Value = 2830.5
$ws.Range("L63").Value = 2633.3333
$ws.Range("M63").Value = -2144.5
$ws.Range("N63").Value = -4005.3333
# Row 66
$ws.Range("H66").Value = 2785
$ws.Range("I66").Value = 2830.5
$ws.Range("J66").Value = 2633.3333
$ws.Range("K66").Value = 14152.5
$ws.Range("L66").Value = 13166.6665
$ws.Range("M66").Value = -10720.5
$ws.Range("N66").Value = -20030.6665
# Row 116
$ws.Range("H116").Value = 1312.6818
$ws.Range("I116").Value = 938.8163500000001
$ws.Range("J116").Value = 2390.2942
$ws.Range("K116").Value = 938.8163500000001
$ws.Range("L116").Value = 2390.2942
$ws.Range("M116").Value = 1355.18365
$ws.Range("N116").Value = -6978.2942
# Row 132
$ws.Range("H132").Value = 35230.445
$ws.Range("I132").Value = 28254.342
$ws.Range("J132").Value = 45834.12
$ws.Range("K132").Value = 84763.026
$ws.Range("L132").Value = 137502.36
$ws.Range("M132").Value = -82233.026
$ws.Range("N132").Value = -142562.36
# Row 139
$ws.Range("H139").Value = 49045.555
$ws.Range("J139").Value = 49045.555
$ws.Range("L139").Value = 49045.555
$ws.Range("N139").Value = -59325.555

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1312.6818
$ws.Range("I3").Value = 938.8163500000001
$ws.Range("J3").Value = 2390.2942
$ws.Range("K3").Value = 938.8163500000001
$ws.Range("L3").Value = 2390.2942
$ws.Range("M3").Value = -824.8163500000001
$ws.Range("N3").Value = -2618.2942
# Row 94
$ws.Range("H94").Value = 673.3333
$ws.Range("I94").Value = 637.8570999999999
$ws.Range("K94").Value = 637.8570999999999
$ws.Range("M94").Value = -186.8570999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 33549.383
$ws.Range("I31").Value = 28398.572
$ws.Range("J31").Value = 36730.766
$ws.Range("K31").Value = 28398.572
$ws.Range("L31").Value = 36730.766
$ws.Range("M31").Value = -28103.572
$ws.Range("N31").Value = -37320.766
# Row 34
$ws.Range("H34").Value = 33549.383
$ws.Range("I34").Value = 28398.572
$ws.Range("J34").Value = 36730.766
$ws.Range("K34").Value = 28398.572
$ws.Range("L34").Value = 36730.766
$ws.Range("M34").Value = -28196.572
$ws.Range("N34").Value = -37134.766

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 653.72095
$ws.Range("I5").Value = 327
$ws.Range("J5").Value = 1330.5
$ws.Range("K5").Value = 981
$ws.Range("L5").Value = 3991.5
$ws.Range("M5").Value = -869
$ws.Range("N5").Value = -4215.5
# Row 34
$ws.Range("H34").Value = 851.4483
$ws.Range("I34").Value = 300
$ws.Range("J34").Value = 1061.5238
$ws.Range("K34").Value = 900
$ws.Range("L34").Value = 3184.5714
$ws.Range("M34").Value = -816
$ws.Range("N34").Value = -3352.5714
# Row 39
$ws.Range("H39").Value = 2115
$ws.Range("J39").Value = 2322.5
$ws.Range("L39").Value = 6967.5
$ws.Range("N39").Value = -7555.5
# Row 56
$ws.Range("H56").Value = 23240.723
$ws.Range("I56").Value = 23240.723
$ws.Range("K56").Value = 23240.723
$ws.Range("M56").Value = -22710.723
# Row 122
$ws.Range("H122").Value = 1084
$ws.Range("I122").Value = 364.8
$ws.Range("J122").Value = 1273.2632
$ws.Range("K122").Value = 3283.2
$ws.Range("L122").Value = 11459.3688
$ws.Range("M122").Value = -833.2000000000003
$ws.Range("N122").Value = -16359.3688
# Row 129
$ws.Range("H129").Value = 2453104
$ws.Range("I129").Value = 1716.5834
$ws.Range("J129").Value = 3790224.2
$ws.Range("K129").Value = 5149.7502
$ws.Range("L129").Value = 11370672.6
$ws.Range("M129").Value = -149.7502000000004
$ws.Range("N129").Value = -11380672.6
# Row 131
$ws.Range("H131").Value = 927.7241
$ws.Range("J131").Value = 1001.72546
$ws.Range("L131").Value = 3005.17638
$ws.Range("N131").Value = -13085.17638
# Row 133
$ws.Range("H133").Value = 6332.36
$ws.Range("I133").Value = 3278.889
$ws.Range("J133").Value = 8049.9375
$ws.Range("K133").Value = 9836.667000000001
$ws.Range("L133").Value = 24149.8125
$ws.Range("M133").Value = -4776.667000000001
$ws.Range("N133").Value = -34269.8125
# Row 135
$ws.Range("H135").Value = 653.72095
$ws.Range("I135").Value = 327
$ws.Range("J135").Value = 1330.5
$ws.Range("K135").Value = 2943
$ws.Range("L135").Value = 11974.5
$ws.Range("M135").Value = -408
$ws.Range("N135").Value = -17044.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1407.9333
$ws.Range("I97").Value = 1354.32
$ws.Range("J97").Value = 1676
$ws.Range("K97").Value = 1354.32
$ws.Range("L97").Value = 1676
$ws.Range("M97").Value = -858.3199999999999
$ws.Range("N97").Value = -2668
# Row 141
$ws.Range("H141").Value = 37240
$ws.Range("J141").Value = 37240
$ws.Range("L141").Value = 37240
$ws.Range("N141").Value = -47600

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 956.4194
$ws.Range("I93").Value = 872
$ws.Range("J93").Value = 1526.25
$ws.Range("K93").Value = 872
$ws.Range("L93").Value = 1526.25
$ws.Range("M93").Value = 376
$ws.Range("N93").Value = -4022.25
# Row 136
$ws.Range("H136").Value = 68525.67999999999
$ws.Range("I136").Value = 40596.15
$ws.Range("J136").Value = 257050
$ws.Range("K136").Value = 121788.45
$ws.Range("L136").Value = 771150
$ws.Range("M136").Value = -119238.45
$ws.Range("N136").Value = -776250

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 140
$ws.Range("H140").Value = 51094.75
$ws.Range("J140").Value = 51094.75
$ws.Range("L140").Value = 51094.75
$ws.Range("N140").Value = -61454.75
